$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 2308
$ws.Cells.Item(7, 6).Value = 8180
$ws.Cells.Item(8, 6).Value = 121
$ws.Cells.Item(10, 6).Value = 1623
$ws.Cells.Item(11, 6).Value = 1331
$ws.Cells.Item(12, 6).Value = 212
$ws.Cells.Item(13, 6).Value = 4516
$ws.Cells.Item(14, 6).Value = 6184
$ws.Cells.Item(15, 6).Value = 797
$ws.Cells.Item(16, 6).Value = 61
$ws.Cells.Item(17, 6).Value = 1260
$ws.Cells.Item(19, 6).Value = 488
$ws.Cells.Item(20, 6).Value = 6530
$ws.Cells.Item(21, 6).Value = 363
$ws.Cells.Item(24, 6).Value = 4401
$ws.Cells.Item(25, 6).Value = 320
$ws.Cells.Item(26, 6).Value = 724
$ws.Cells.Item(27, 6).Value = 2058
$ws.Cells.Item(28, 6).Value = 1197
$ws.Cells.Item(29, 6).Value = 358
$ws.Cells.Item(31, 6).Value = 72
$ws.Cells.Item(33, 6).Value = 46
$ws.Cells.Item(34, 6).Value = 91
$ws.Cells.Item(35, 6).Value = 336
$ws.Cells.Item(36, 6).Value = 1202
$ws.Cells.Item(37, 6).Value = 1915
$ws.Cells.Item(38, 6).Value = 152
$ws.Cells.Item(39, 6).Value = 446
$ws.Cells.Item(40, 6).Value = 171
$ws.Cells.Item(41, 6).Value = 1232
$ws.Cells.Item(44, 6).Value = 1197
$ws.Cells.Item(47, 6).Value = 201
$ws.Cells.Item(48, 6).Value = 34
$ws.Cells.Item(49, 6).Value = 24

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(7, 6).Value = 404
$ws.Cells.Item(8, 6).Value = 419
$ws.Cells.Item(10, 6).Value = 220
$ws.Cells.Item(17, 6).Value = 109
$ws.Cells.Item(22, 6).Value = 129
$ws.Cells.Item(26, 6).Value = 178

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(8, 6).Value = 3164
$ws.Cells.Item(9, 6).Value = 1086
$ws.Cells.Item(10, 6).Value = 1165
$ws.Cells.Item(11, 6).Value = 1526
$ws.Cells.Item(12, 6).Value = 1867
$ws.Cells.Item(13, 6).Value = 356
$ws.Cells.Item(14, 6).Value = 223

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 6).Value = 3164
$ws.Cells.Item(8, 6).Value = 2308
$ws.Cells.Item(9, 6).Value = 1086
$ws.Cells.Item(11, 6).Value = 1623
$ws.Cells.Item(12, 6).Value = 1526
$ws.Cells.Item(13, 6).Value = 1331
$ws.Cells.Item(15, 6).Value = 212
$ws.Cells.Item(16, 6).Value = 1867
$ws.Cells.Item(17, 6).Value = 4516
$ws.Cells.Item(18, 6).Value = 356
$ws.Cells.Item(19, 6).Value = 404
$ws.Cells.Item(20, 6).Value = 419
$ws.Cells.Item(21, 6).Value = 797
$ws.Cells.Item(22, 6).Value = 61
$ws.Cells.Item(23, 6).Value = 1260
$ws.Cells.Item(25, 6).Value = 488
$ws.Cells.Item(26, 6).Value = 6530
$ws.Cells.Item(27, 6).Value = 363
$ws.Cells.Item(28, 6).Value = 223
$ws.Cells.Item(30, 6).Value = 4401
$ws.Cells.Item(31, 6).Value = 320
$ws.Cells.Item(32, 6).Value = 2058
$ws.Cells.Item(33, 6).Value = 1197
$ws.Cells.Item(34, 6).Value = 358
$ws.Cells.Item(36, 6).Value = 72
$ws.Cells.Item(37, 6).Value = 46
$ws.Cells.Item(39, 6).Value = 91
$ws.Cells.Item(40, 6).Value = 336
$ws.Cells.Item(41, 6).Value = 1915
$ws.Cells.Item(42, 6).Value = 152
$ws.Cells.Item(43, 6).Value = 446
$ws.Cells.Item(44, 6).Value = 1232
$ws.Cells.Item(45, 6).Value = 129
$ws.Cells.Item(48, 6).Value = 1197
$ws.Cells.Item(49, 6).Value = 201
